$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 updates ---
$ws.Range("B2").Value = 3
$ws.Range("C2").Value = "some code"
$ws.Range("BA2").Value = 25
$ws.Range("CH2").Value = "WFR"

# --- Row 3 updates ---
$ws.Range("B3").Value = 4
$ws.Range("C3").Value = "dasda"
$ws.Range("N3").Value = ""
$ws.Range("BA3").Value = 3
$ws.Range("CF3").Value = "test1"
$ws.Range("CG3").Value = "PCC"
$ws.Range("CH3").Value = "REF"

# --- Row 4 new row ---
# Copy A3's format (bold/border style used for column A) onto A4, then set its value.
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 5
$ws.Range("C4").Value = "sadasd"
$ws.Range("BA4").Value = 8
$ws.Range("CF4").Value = "central store"
$ws.Range("CG4").Value = "TEM"
$ws.Range("CH4").Value = "WFR"
